# Updated status for UDQ pages
# - Several "In Progress" / "Not Done" statuses in the UDQ table (Sheet1, rows 5-38)
#   are updated to "Done".
# - The corresponding row-label cells in column B lose their special bold/colored
#   "highlight" font, reverting to the plain look used by the rest of the table
#   (they keep their yellow background highlight).
# - The now-unused "In Progress" shared string / styles are naturally dropped by
#   Excel once nothing references them any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# --- Reference cells that already carry the formatting we want to re-use ---
$doneRef  = $ws.Range("C5")   # plain "Done" look: no fill, green font
$plain7   = $ws.Range("B17")  # plain label look, style used by most rows (odd)
$plain8   = $ws.Range("B6")   # plain label look, style used by most rows (even)

# --- 1. Flip stale statuses to "Done" ------------------------------------
$statusCells = @("C31","D7","E7","D31","E31","D33","E33","D34","E34","D35","E35","D37","E37")

foreach ($addr in $statusCells) {
    $cell = $ws.Range($addr)
    $doneRef.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $cell.Value = "Done"
}

# --- 2. Normalize the highlighted row-label cells -------------------------
$toPlain7 = @("B18","B20","B22","B23","B24","B26","B28","B30")
foreach ($addr in $toPlain7) {
    $cell = $ws.Range($addr)
    $plain7.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
}

$plain8.Copy() | Out-Null
$ws.Range("B38").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0
